$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2410" and "_new" -> "_FV2504" ---
# Columns A-J: "_old" suffixed headers; column K: "diff" (unchanged);
# Columns L-U: "_new" suffixed headers.
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- 2. Turn the used range into a real Excel Table ("Table1") ---
# Stash the header row's existing formatting (fill/border/font/alignment) in
# a scratch range so it survives ListObjects.Add() (which otherwise bakes
# the header row's current look into a brand-new dxf/table style).
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")
$headerRange.Copy()
$scratch.PasteSpecial(-4122) # xlPasteFormats
$headerRange.ClearFormats()

$list = $ws.ListObjects.Add(1, $ws.Range("A1:U91"), $null, 1)
$list.Name = "Table1"
$list.TableStyle = ""

# Restore the header row's original formatting and clean up the scratch area.
$scratch.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats
$scratch.Clear()
$excel.CutCopyMode = $false

# --- 3. Freeze the header row (row 1) via a split pane ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
